$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the "ivide, " / "multiplication,..." / ")" runs
#           (dropping the proofErr markers) into a single run, and
#           split "Expressions (if else, switch)" into:
#             "Expressions"
#             "Logic " + "(if else, switch)"
# -----------------------------------------------------------------
$pComputations = $null
$pExpressions = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Computations (plus, minus, modulus*") {
        $pComputations = $p
    }
    if ($t -like "Expressions (if else, switch)*") {
        $pExpressions = $p
    }
}

$rng1 = $d.Range($pComputations.Range.Start, $pExpressions.Range.End)
$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>Computations (plus, minus, modulus, d</w:t></w:r><w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>ivide, multiplication,&#8230;)</w:t></w:r></w:p>
<w:p><w:r><w:t>Expressions</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Logic </w:t></w:r><w:r><w:t>(if else, switch)</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng1.InsertXML($xml1) | Out-Null

# -----------------------------------------------------------------
# Change 2: move <w:lastRenderedPageBreak/> from before "Pointers"
#           to before "Functions" (Week 5 section)
# -----------------------------------------------------------------
$pFunctions = $null
$pPointers = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Functions`r") {
        $pFunctions = $p
    }
    if ($t -eq "Pointers`r") {
        $pPointers = $p
    }
}

$rng2 = $d.Range($pFunctions.Range.Start, $pPointers.Range.End)
$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r><w:lastRenderedPageBreak/><w:t>Functions</w:t></w:r></w:p>
<w:p><w:r><w:t>Pointers</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng2.InsertXML($xml2) | Out-Null

# -----------------------------------------------------------------
# Change 3: move <w:lastRenderedPageBreak/> from before
#           "Language Standards" to before "Standard Library"
# -----------------------------------------------------------------
$pStdLib = $null
$pLangStd = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Standard Library*") {
        $pStdLib = $p
    }
    if ($t -eq "Language Standards`r") {
        $pLangStd = $p
    }
}

$rng3 = $d.Range($pStdLib.Range.Start, $pLangStd.Range.End)
$xml3 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r><w:lastRenderedPageBreak/><w:t>Standard Library</w:t></w:r><w:r><w:t xml:space="preserve"> (iostream)</w:t></w:r></w:p>
<w:p><w:r><w:t>Language Standards</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng3.InsertXML($xml3) | Out-Null

Write-Host "Edits applied"
